$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '25.078.92'
Set-TextValue $ws.Range("E2") '  -3.03%  '

Set-TextValue $ws.Range("D3") '1.648.98'
Set-TextValue $ws.Range("E3") '  -5.02%  '

Set-TextValue $ws.Range("D4") '0.9999'
Set-TextValue $ws.Range("E4") '  +0.01%  '

Set-TextValue $ws.Range("D5") '237.06'
Set-TextValue $ws.Range("E5") '  -2.20%  '

Set-TextValue $ws.Range("E6") '  +0.00%  '

Set-TextValue $ws.Range("D7") '0.4786'
Set-TextValue $ws.Range("E7") '  -7.98%  '

Set-TextValue $ws.Range("D8") '0.2622'
Set-TextValue $ws.Range("E8") '  -4.35%  '

Set-TextValue $ws.Range("D9") '0.06033'
Set-TextValue $ws.Range("E9") '  -2.05%  '

Set-TextValue $ws.Range("E10") '  -1.17%  '

Set-TextValue $ws.Range("D11") '1.647.44'
Set-TextValue $ws.Range("E11") '  -5.15%  '

Set-TextValue $ws.Range("D12") '14.46'
Set-TextValue $ws.Range("E12") '  -3.43%  '

Set-TextValue $ws.Range("D13") '0.6180'
Set-TextValue $ws.Range("E13") '  -3.92%  '

Set-TextValue $ws.Range("D14") '4.565'
Set-TextValue $ws.Range("E14") '  -1.03%  '

Set-TextValue $ws.Range("D15") '73.06'
Set-TextValue $ws.Range("E15") '  -5.48%  '

Set-TextValue $ws.Range("D16") '0.9997'
Set-TextValue $ws.Range("E16") '  -0.03%  '

Set-TextValue $ws.Range("D17") '0.9999'
Set-TextValue $ws.Range("E17") '  +0.00%  '

Set-TextValue $ws.Range("D18") '25.062.66'
Set-TextValue $ws.Range("E18") '  -3.22%  '

Set-TextValue $ws.Range("D20") '0.000006563'
Set-TextValue $ws.Range("E20") '  -3.18%  '

Set-TextValue $ws.Range("D21") '4.421'
Set-TextValue $ws.Range("E21") '  +3.21%  '

Set-TextValue $ws.Range("D22") '1.865.20'
Set-TextValue $ws.Range("E22") '  -4.97%  '

Set-TextValue $ws.Range("D23") '8.469'

Set-TextValue $ws.Range("D24") '5.235'
Set-TextValue $ws.Range("E24") '  -1.08%  '

Set-TextValue $ws.Range("D25") '133.86'
Set-TextValue $ws.Range("E25") '  -2.42%  '

Set-TextValue $ws.Range("D26") '14.71'
Set-TextValue $ws.Range("E26") '  -3.36%  '

Set-TextValue $ws.Range("D27") '1.396'
Set-TextValue $ws.Range("E27") '  -7.55%  '

Set-TextValue $ws.Range("D28") '1.689'
Set-TextValue $ws.Range("E28") '  -4.75%  '

Set-TextValue $ws.Range("D29") '101.84'
Set-TextValue $ws.Range("E29") '  -3.20%  '

Set-TextValue $ws.Range("E30") '  -4.50%  '

Set-TextValue $ws.Range("D31") '0.07926'
Set-TextValue $ws.Range("E31") '  -3.98%  '

Set-TextValue $ws.Range("D32") '3.548'
Set-TextValue $ws.Range("E32") '  -2.84%  '

Set-TextValue $ws.Range("D33") '0.04546'
Set-TextValue $ws.Range("E33") '  -2.64%  '

Set-TextValue $ws.Range("D34") '2.607'
Set-TextValue $ws.Range("E34") '  -2.04%  '

Set-TextValue $ws.Range("D35") '0.9397'
Set-TextValue $ws.Range("E35") '  -5.07%  '

Set-TextValue $ws.Range("D36") '0.5785'
Set-TextValue $ws.Range("E36") '  -6.73%  '

Set-TextValue $ws.Range("D37") '2.625'
Set-TextValue $ws.Range("E37") '  -2.34%  '

Set-TextValue $ws.Range("D38") '0.01537'
Set-TextValue $ws.Range("E38") '  -3.92%  '

Set-TextValue $ws.Range("B39") 'TrustWalletToken'
Set-TextValue $ws.Range("C39") 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D39") '0.8394'
Set-TextValue $ws.Range("E39") '  +12.40%  '

Set-TextValue $ws.Range("B40") 'PaxDollar'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range("D40") '0.9997'
Set-TextValue $ws.Range("E40") '  +0.01%  '

Set-TextValue $ws.Range("D41") '1.818'
Set-TextValue $ws.Range("E41") '  -5.37%  '

Set-TextValue $ws.Range("D42") '98.67'
Set-TextValue $ws.Range("E42") '  -1.17%  '

Set-TextValue $ws.Range("D43") '0.3709'
Set-TextValue $ws.Range("E43") '  -3.88%  '

Set-TextValue $ws.Range("D44") '4.796'
Set-TextValue $ws.Range("E44") '  -4.11%  '

Set-TextValue $ws.Range("D45") '0.1121'

Set-TextValue $ws.Range("D46") '6.039'
Set-TextValue $ws.Range("E46") '  -3.47%  '

Set-TextValue $ws.Range("D47") '0.05156'
Set-TextValue $ws.Range("E47") '  -1.07%  '

Set-TextValue $ws.Range("D48") '52.15'
Set-TextValue $ws.Range("E48") '  -5.27%  '

Set-TextValue $ws.Range("D49") '29.44'
Set-TextValue $ws.Range("E49") '  -3.81%  '

Set-TextValue $ws.Range("D50") '1.000'
Set-TextValue $ws.Range("E50") '  -0.02%  '

Set-TextValue $ws.Range("D51") '0.3332'
Set-TextValue $ws.Range("E51") '  -2.37%  '
